$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.725.11"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "2.385.97"
$ws.Range("E3").Value = "  -4.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.78"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.94"
$ws.Range("E6").Value = "  -4.31%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  -8.09%  "
$ws.Range("D9").Value = "2.387.30"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.338"
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.13"
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "2.817.89"
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("D16").Value = "60.725.56"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000161"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "2.386.26"
$ws.Range("E18").Value = "  -4.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("E19").Value = "  -5.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.66"
$ws.Range("E20").Value = "  -5.23%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.02"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "310.70"
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.77"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.85"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "2.513.54"
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("D28").Value = "0.0₃0913"
$ws.Range("E28").Value = "  -10.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.40"
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -5.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "504.23"
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.143"
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.40"
$ws.Range("E37").Value = "  -8.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.54"
$ws.Range("E38").Value = "  -7.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.82"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "138.07"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("E45").Value = "  -7.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "138.43"
$ws.Range("E46").Value = "  -7.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.47"
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.95"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  -4.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.573"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0912"
$ws.Range("E51").Value = "  -4.39%  "
